$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right marks 4 -> 5, Wrong marks -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total"): Right total 40 -> 50, Wrong total -4 -> -4.8,
# and the "obtained/max" summary string 36/112 -> 45.2/140
$ws.Range("B12").Value = 50
$ws.Range("C12").Value = -4.8
$ws.Range("E12").Value = "45.2/140"
